$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.406.31'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '1.872.69'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '245.57'
$ws.Range('E5').Value = '  -0.93%  '
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = '0.4713'
$ws.Range('E7').Value = '  -0.56%  '
$ws.Range('D8').Value = '0.2866'
$ws.Range('E8').Value = '  -1.98%  '
$ws.Range('D9').Value = '0.06495'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('D10').Value = '21.86'
$ws.Range('E10').Value = '  -0.51%  '
$ws.Range('D11').Value = '100.44'
$ws.Range('E11').Value = '  +3.56%  '
$ws.Range('D12').Value = '0.07806'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').Value = '1.870.20'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').Value = '0.7287'
$ws.Range('E14').Value = '  -1.02%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.170'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.48%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '283.90'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').Value = '30.385.36'
$ws.Range('E17').Value = '  -1.32%  '
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.000'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '0.000007485'
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('D21').Value = '2.116.18'
$ws.Range('E21').Value = '  -1.16%  '
$ws.Range('D22').Value = '5.328'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '0.9997'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = '6.335'
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').Value = '9.045'
$ws.Range('E25').Value = '  -1.88%  '
$ws.Range('D26').Value = '161.99'
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('D27').Value = '18.96'
$ws.Range('E27').Value = '  +0.26%  '
$ws.Range('D28').Value = '1.896'
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('D29').Value = '0.09686'
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('E30').Value = '  -1.47%  '
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').Value = '4.226'
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('D33').Value = '4.154'
$ws.Range('E33').Value = '  -0.98%  '
$ws.Range('D34').Value = '0.04811'
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('D35').Value = '1.125'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6900'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').Value = '2.708'
$ws.Range('E37').Value = '  -0.55%  '
$ws.Range('D38').Value = '0.01894'
$ws.Range('D39').Value = '2.843'
$ws.Range('E39').Value = '  +1.14%  '
$ws.Range('D40').Value = '76.11'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('D41').Value = '6.305'
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D42').Value = '1.953'
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('D43').Value = '0.4215'
$ws.Range('E43').Value = '  -0.95%  '
$ws.Range('D44').Value = '0.9993'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').Value = '0.8248'
$ws.Range('E45').Value = '  -1.29%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '100.80'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('D47').Value = '9.728'
$ws.Range('E47').Value = '  +2.29%  '
$ws.Range('D48').Value = '7.012'
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('E49').Value = '  -1.78%  '
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('D51').Value = '883.93'
$ws.Range('E51').Value = '  -3.55%  '
